$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.718.33"
$ws.Range("E2").Value = "  +10.93%  "

$ws.Range("D3").Value = "3.379.28"
$ws.Range("E3").Value = "  +7.82%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "651.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.73%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.402"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +44.56%  "

$ws.Range("E8").Value = "  -0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.57%  "

$ws.Range("D10").Value = "3.375.19"
$ws.Range("E10").Value = "  +7.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.596"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.64%  "

$ws.Range("E12").Value = "  +15.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +18.90%  "

$ws.Range("E14").Value = "  +2.65%  "

$ws.Range("D15").Value = "3.966.74"
$ws.Range("E15").Value = "  +6.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.96%  "

$ws.Range("D17").Value = "89.318.41"
$ws.Range("E17").Value = "  +10.62%  "

$ws.Range("D18").Value = "3.340.32"
$ws.Range("E18").Value = "  +6.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "460.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.14%  "

$ws.Range("E24").Value = "  +4.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +21.69%  "

$ws.Range("D27").Value = "3.480.72"
$ws.Range("E27").Value = "  +5.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000143"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +19.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "78.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.66%  "

$ws.Range("E30").Value = "  +45.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "599.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.87%  "

$ws.Range("E35").Value = "  -2.42%  "

$ws.Range("E36").Value = "  +9.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +26.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.144"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.428"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.65%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.95%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +13.07%  "

$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "157.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "190.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.46%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.677"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.69%  "
